$d = $word.ActiveDocument

# 1. Merge the runs describing the 5 files / classes sentence into a single run.
#    (Text content is unchanged; this just collapses several split runs, and the
#    gramStart/gramEnd proofErr markers around "class", into one run.)
$d.Content.Find.Execute(
    "My code consisted of 5 files, a main file, and airport class, an airline class, a route class and a node class. Each class was split into a header and cpp file, the header containing the basic member variables and function declarations and the cpp containing the implementations of the various functions each class had.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "My code consisted of 5 files, a main file, and airport class, an airline class, a route class and a node class. Each class was split into a header and cpp file, the header containing the basic member variables and function declarations and the cpp containing the implementations of the various functions each class had.",
    2
) | Out-Null

# 2. Merge the runs describing the airport class attributes into a single run.
$d.Content.Find.Execute(
    "The airport class contained 6 attributes from the dataset such as the airport id, airport name etc. The Airline class has 8 attributes such as the Airline ID, the name, the alias/nickname of the airline or the IATA or ICAO codes. The route class contains 9 attributes, such as the source and destination airline codes, the costs, and the airline Id to name a few. For each class, a ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The airport class contained 6 attributes from the dataset such as the airport id, airport name etc. The Airline class has 8 attributes such as the Airline ID, the name, the alias/nickname of the airline or the IATA or ICAO codes. The route class contains 9 attributes, such as the source and destination airline codes, the costs, and the airline Id to name a few. For each class, a ",
    2
) | Out-Null

# 3. Merge the runs describing the Retrieve_ method into a single run.
$d.Content.Find.Execute(
    "method was created whose goal was to read and extract the relevant data from the respective csv files and store them in an unordered map",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "method was created whose goal was to read and extract the relevant data from the respective csv files and store them in an unordered map",
    2
) | Out-Null

# 4. Append the new sentence about the breadth first search algorithm.
$rng = $d.Content
$rng.Find.Execute(
    "Due to the lack of certain key data structures and functions in C++, I struggled quite a bit in completing this project and implementing some of the algorithms such as the",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "",
    0
) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(" breadth first search algorithm and finding the appropriate data structures to store them in.")
